$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 41
$ws.Range("H41").Value2 = 183.29411
$ws.Range("I41").Value2 = 208.88889
$ws.Range("J41").Value2 = 154.5
$ws.Range("K41").Value2 = 208.88889
$ws.Range("L41").Value2 = 154.5
$ws.Range("M41").Value2 = 231.11111
$ws.Range("N41").Value2 = -1034.5
# Row 76
$ws.Range("H76").Value2 = 3320.2307
$ws.Range("I76").Value2 = 3283
$ws.Range("J76").Value2 = 3444.3333
$ws.Range("K76").Value2 = 3283
$ws.Range("L76").Value2 = 3444.3333
$ws.Range("M76").Value2 = -2968
$ws.Range("N76").Value2 = -4074.3333
# Row 79
$ws.Range("H79").Value2 = 3320.2307
$ws.Range("I79").Value2 = 3283
$ws.Range("J79").Value2 = 3444.3333
$ws.Range("K79").Value2 = 3283
$ws.Range("L79").Value2 = 3444.3333
$ws.Range("M79").Value2 = -2191
$ws.Range("N79").Value2 = -5628.3333
# Row 86
$ws.Range("H86").Value2 = 4719.7
$ws.Range("I86").Value2 = 1998
$ws.Range("J86").Value2 = 5022.1113
$ws.Range("K86").Value2 = 1998
$ws.Range("L86").Value2 = 5022.1113
$ws.Range("M86").Value2 = -875
$ws.Range("N86").Value2 = -7268.1113
# Row 89
$ws.Range("H89").Value2 = 4719.7
$ws.Range("I89").Value2 = 1998
$ws.Range("J89").Value2 = 5022.1113
$ws.Range("K89").Value2 = 9990
$ws.Range("L89").Value2 = 25110.5565
$ws.Range("M89").Value2 = -4374
$ws.Range("N89").Value2 = -36342.5565
# Row 106
$ws.Range("H106").Value2 = 987.6875
$ws.Range("I106").Value2 = 993.5294
$ws.Range("J106").Value2 = 981.06665
$ws.Range("K106").Value2 = 993.5294
$ws.Range("L106").Value2 = 981.06665
$ws.Range("M106").Value2 = -362.5294
$ws.Range("N106").Value2 = -2243.06665
# Row 112
$ws.Range("H112").Value2 = 47620428
$ws.Range("I112").Value2 = 566.6667
$ws.Range("J112").Value2 = 55557070
$ws.Range("K112").Value2 = 1700.0001
$ws.Range("L112").Value2 = 166671210
$ws.Range("M112").Value2 = -592.0001
$ws.Range("N112").Value2 = -166673426

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 37
$ws.Range("H37").Value2 = 8572
$ws.Range("J37").Value2 = 10587.714
$ws.Range("L37").Value2 = 10587.714
$ws.Range("N37").Value2 = -11133.714
# Row 45
$ws.Range("H45").Value2 = 1149.1578
$ws.Range("I45").Value2 = 934.8570999999999
$ws.Range("J45").Value2 = 1274.1666
$ws.Range("K45").Value2 = 934.8570999999999
$ws.Range("L45").Value2 = 1274.1666
$ws.Range("M45").Value2 = -557.8570999999999
$ws.Range("N45").Value2 = -2028.1666
# Row 138
$ws.Range("H138").Value2 = 56328
$ws.Range("J138").Value2 = 56328
$ws.Range("L138").Value2 = 56328
$ws.Range("N138").Value2 = -66608
# Row 139
$ws.Range("H139").Value2 = 56891.668
$ws.Range("J139").Value2 = 56891.668
$ws.Range("L139").Value2 = 56891.668
$ws.Range("N139").Value2 = -67171.66800000001
# Row 140
$ws.Range("H140").Value2 = 56959.6
$ws.Range("J140").Value2 = 56959.6
$ws.Range("L140").Value2 = 56959.6
$ws.Range("N140").Value2 = -67319.60000000001
# Row 141
$ws.Range("H141").Value2 = 60632.445
$ws.Range("J141").Value2 = 60632.445
$ws.Range("L141").Value2 = 60632.445
$ws.Range("N141").Value2 = -70992.44500000001

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value2 = 62807.58
$ws.Range("I134").Value2 = 74147.19
$ws.Range("J134").Value2 = 2329.6667
$ws.Range("K134").Value2 = 222441.57
$ws.Range("L134").Value2 = 6989.000100000001
$ws.Range("M134").Value2 = -219906.57
$ws.Range("N134").Value2 = -12059.0001
# Row 138
$ws.Range("H138").Value2 = 61532.715
$ws.Range("J138").Value2 = 61532.715
$ws.Range("L138").Value2 = 61532.715
$ws.Range("N138").Value2 = -71812.715
# Row 140
$ws.Range("H140").Value2 = 58979
$ws.Range("J140").Value2 = 58979
$ws.Range("L140").Value2 = 58979
$ws.Range("N140").Value2 = -69339

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 59
$ws.Range("H59").Value2 = 15499.875
$ws.Range("J59").Value2 = 15428.429
$ws.Range("L59").Value2 = 15428.429
$ws.Range("N59").Value2 = -17718.429
# Row 135
$ws.Range("H135").Value2 = 58964.445
$ws.Range("J135").Value2 = 58964.445
$ws.Range("L135").Value2 = 58964.445
$ws.Range("N135").Value2 = -69104.44500000001
# Row 138
$ws.Range("H138").Value2 = 63147.273
$ws.Range("I138").Value2 = 9800
$ws.Range("J138").Value2 = 68482
$ws.Range("K138").Value2 = 9800
$ws.Range("L138").Value2 = 68482
$ws.Range("M138").Value2 = -4660
$ws.Range("N138").Value2 = -78762
# Row 140
$ws.Range("H140").Value2 = 150371.11
$ws.Range("J140").Value2 = 179048.58
$ws.Range("L140").Value2 = 179048.58
$ws.Range("N140").Value2 = -189408.58
# Row 141
$ws.Range("H141").Value2 = 48320.312
$ws.Range("J141").Value2 = 48320.312
$ws.Range("L141").Value2 = 48320.312
$ws.Range("N141").Value2 = -58680.312

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 117
$ws.Range("H117").Value2 = 1146.4445
$ws.Range("I117").Value2 = 932.6667
$ws.Range("J117").Value2 = 1574
$ws.Range("K117").Value2 = 2798.0001
$ws.Range("L117").Value2 = 4722
$ws.Range("M117").Value2 = 643.9998999999998
$ws.Range("N117").Value2 = -11606
# Row 122
$ws.Range("H122").Value2 = 7408414.5
$ws.Range("I122").Value2 = 10417197
$ws.Range("J122").Value2 = 2180
$ws.Range("K122").Value2 = 93754773
$ws.Range("L122").Value2 = 19620
$ws.Range("M122").Value2 = -93752323
$ws.Range("N122").Value2 = -24520
# Row 129
$ws.Range("H129").Value2 = 1632.8572
$ws.Range("I129").Value2 = 743
$ws.Range("J129").Value2 = 1781.1666
$ws.Range("K129").Value2 = 2229
$ws.Range("L129").Value2 = 5343.4998
$ws.Range("M129").Value2 = 2771
$ws.Range("N129").Value2 = -15343.4998
# Row 131
$ws.Range("H131").Value2 = 2468.6155
$ws.Range("I131").Value2 = 20110
$ws.Range("J131").Value2 = 1615
$ws.Range("K131").Value2 = 60330
$ws.Range("L131").Value2 = 4845
$ws.Range("M131").Value2 = -55290
$ws.Range("N131").Value2 = -14925

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 136
$ws.Range("H136").Value2 = 1838.7142
$ws.Range("I136").Value2 = 1083.1177
$ws.Range("K136").Value2 = 3249.3531
$ws.Range("M136").Value2 = -699.3531000000003
# Row 139
$ws.Range("H139").Value2 = 53156.363
$ws.Range("J139").Value2 = 53156.363
$ws.Range("L139").Value2 = 53156.363
$ws.Range("N139").Value2 = -63436.363
# Row 141
$ws.Range("H141").Value2 = 66560
$ws.Range("I141").Value2 = 8800
$ws.Range("J141").Value2 = 81000
$ws.Range("K141").Value2 = 8800
$ws.Range("L141").Value2 = 81000
$ws.Range("M141").Value2 = -3620
$ws.Range("N141").Value2 = -91360

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 136
$ws.Range("H136").Value2 = 3578807
$ws.Range("I136").Value2 = 11503.479
$ws.Range("J136").Value2 = 8405159
$ws.Range("K136").Value2 = 34510.437
$ws.Range("L136").Value2 = 25215477
$ws.Range("M136").Value2 = -31960.437
$ws.Range("N136").Value2 = -25220577
